$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column U: header "06-10-2020" in U1, mirroring the formatting of
# the preceding date header cells (bold, thin border, centered/top aligned).
# xlCenter = -4108, xlTop = -4160, xlContinuous = 1
$hdr = $ws.Range("U1")
$hdr.Value = "placeholder"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.NumberFormat = "@"
$hdr.Value = "06-10-2020"

# --- New column U values for each state/UT row (U2:U36) ---
$values = @(
    3659,
    666433,
    7775,
    153491,
    176995,
    10797,
    97067,
    2991,
    263938,
    30456,
    123638,
    121596,
    12653,
    63790,
    76843,
    522846,
    149111,
    3414,
    115878,
    1162585,
    9334,
    4491,
    1837,
    5422,
    206400,
    24221,
    102648,
    123421,
    2547,
    569664,
    174769,
    22131,
    42621,
    366321,
    240707
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 21).Value = $values[$i]
}
